$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells whose new values look numeric stay as text (matches original inlineStr formatting)
$textCells = @("D5","D9","D10","D11","D14","D17","D20","D23","D24","D25","D26","D27","D28","D29","D31","D32","D35","D41","D42","D43","D47")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values from the data refresh
$ws.Range("D2").Value = "25.276.91"
$ws.Range("E2").Value = "  -2.50%  "
$ws.Range("D3").Value = "1.562.18"
$ws.Range("E3").Value = "  -3.78%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "206.84"
$ws.Range("E5").Value = "  -3.18%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  -4.95%  "
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("D9").Value = "0.241"
$ws.Range("E9").Value = "  -3.23%  "
$ws.Range("D10").Value = "17.75"
$ws.Range("E10").Value = "  -2.76%  "
$ws.Range("D11").Value = "0.0782"
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("D12").Value = "1.781.62"
$ws.Range("E12").Value = "  -3.64%  "
$ws.Range("D13").Value = "1.565.26"
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").Value = "3.99"
$ws.Range("E14").Value = "  -4.47%  "
$ws.Range("E15").Value = "  -3.78%  "
$ws.Range("D16").Value = "25.267.36"
$ws.Range("E16").Value = "  -2.50%  "
$ws.Range("D17").Value = "59.14"
$ws.Range("E17").Value = "  -3.29%  "
$ws.Range("E18").Value = "  -3.31%  "
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "184.88"
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("E21").Value = "  -2.68%  "
$ws.Range("E22").Value = "  -3.23%  "
$ws.Range("D23").Value = "5.86"
$ws.Range("E23").Value = "  -3.26%  "
$ws.Range("D24").Value = "1.01"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").Value = "0.127"
$ws.Range("E25").Value = "  -4.02%  "
$ws.Range("D26").Value = "139.47"
$ws.Range("E26").Value = "  -2.93%  "
$ws.Range("D27").Value = "1.63"
$ws.Range("E27").Value = "  -6.87%  "
$ws.Range("D28").Value = "6.48"
$ws.Range("E28").Value = "  -3.41%  "
$ws.Range("D29").Value = "14.79"
$ws.Range("E29").Value = "  -2.29%  "
$ws.Range("E30").Value = "  -6.11%  "
$ws.Range("D31").Value = "0.0463"
$ws.Range("E31").Value = "  -3.96%  "
$ws.Range("D32").Value = "3.02"
$ws.Range("E32").Value = "  -3.10%  "
$ws.Range("E33").Value = "  -4.07%  "
$ws.Range("E34").Value = "  -2.39%  "
$ws.Range("D35").Value = "2.30"
$ws.Range("E35").Value = "  -3.98%  "
$ws.Range("D36").Value = "1.086.76"
$ws.Range("E36").Value = "  -2.88%  "
$ws.Range("E37").Value = "  -0.59%  "
$ws.Range("E38").Value = "  -4.74%  "
$ws.Range("E39").Value = "  -2.18%  "
$ws.Range("E40").Value = "  -5.18%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "0.805"
$ws.Range("E41").Value = "  +4.84%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "0.758"
$ws.Range("E42").Value = "  -10.33%  "
$ws.Range("D43").Value = "92.96"
$ws.Range("E43").Value = "  -4.98%  "
$ws.Range("D45").Value = "1.695.83"
$ws.Range("E45").Value = "  -3.63%  "
$ws.Range("D46").Value = "0.0₆0111"
$ws.Range("E46").Value = "  -2.61%  "
$ws.Range("D47").Value = "52.40"
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("E48").Value = "  -4.71%  "
$ws.Range("E49").Value = "  -2.61%  "
$ws.Range("E50").Value = "  -1.60%  "
$ws.Range("E51").Value = "  -0.43%  "
